# ---------------------------------------------------------------------------
# TP3/Results.xlsx -- "Single classifications tests are done!"
#
# Fills in the G:M (Especificidade..Dados Invalidos) results columns for
# rows 157:187 of the "Single Classification" sheet (previously only A:F
# were populated), and updates both sheets' stored selection/scroll state
# to reflect the newly-filled range.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Single Classification")

# Copy formatting (style) from row 156 G:M block down across rows 157:187
$ws.Range("G156:M156").Copy($ws.Range("G157:M187"))

# Row 157
$ws.Range("G157").Value = 0.90908
$ws.Range("H157").Value = 0.68407
$ws.Range("I157").Value = 747
$ws.Range("J157").Value = 103280
$ws.Range("K157").Value = 10330
$ws.Range("L157").Value = 345
$ws.Range("M157").Value = 0

# Row 158
$ws.Range("G158").Value = 0.17118
$ws.Range("H158").Value = 0.78733
$ws.Range("I158").Value = 733
$ws.Range("J158").Value = 18273
$ws.Range("K158").Value = 88474
$ws.Range("L158").Value = 198
$ws.Range("M158").Value = 7026

# Row 159
$ws.Range("G159").Value = 0.094781
$ws.Range("H159").Value = 0.90221
$ws.Range("I159").Value = 941
$ws.Range("J159").Value = 10438
$ws.Range("K159").Value = 99689
$ws.Range("L159").Value = 102
$ws.Range("M159").Value = 3534

# Row 160
$ws.Range("G160").Value = 0.12901
$ws.Range("H160").Value = 0.80971
$ws.Range("I160").Value = 851
$ws.Range("J160").Value = 14477
$ws.Range("K160").Value = 97739
$ws.Range("L160").Value = 200
$ws.Range("M160").Value = 1437

# Row 161
$ws.Range("G161").Value = 0.95848
$ws.Range("H161").Value = 0.44645
$ws.Range("I161").Value = 371
$ws.Range("J161").Value = 61825
$ws.Range("K161").Value = 2678
$ws.Range("L161").Value = 460
$ws.Range("M161").Value = 49370

# Row 162
$ws.Range("G162").Value = 0.83257
$ws.Range("H162").Value = 0.68541
$ws.Range("I162").Value = 719
$ws.Range("J162").Value = 84592
$ws.Range("K162").Value = 17011
$ws.Range("L162").Value = 330
$ws.Range("M162").Value = 12052

# Row 163
$ws.Range("G163").Value = 0.12469
$ws.Range("H163").Value = 0.79412
$ws.Range("I163").Value = 810
$ws.Range("J163").Value = 13924
$ws.Range("K163").Value = 97745
$ws.Range("L163").Value = 210
$ws.Range("M163").Value = 2015

# Row 164
$ws.Range("G164").Value = 0.9353
$ws.Range("H164").Value = 0.1902
$ws.Range("I164").Value = 194
$ws.Range("J164").Value = 104660
$ws.Range("K164").Value = 7240
$ws.Range("L164").Value = 826
$ws.Range("M164").Value = 1782

# Row 165
$ws.Range("G165").Value = 1
$ws.Range("H165").Value = 0
$ws.Range("I165").Value = 0
$ws.Range("J165").Value = 112000
$ws.Range("K165").Value = 0
$ws.Range("L165").Value = 1047
$ws.Range("M165").Value = 1659

# Row 166
$ws.Range("G166").Value = 0.87379
$ws.Range("H166").Value = 0.11864
$ws.Range("I166").Value = 7
$ws.Range("J166").Value = 90
$ws.Range("K166").Value = 13
$ws.Range("L166").Value = 52
$ws.Range("M166").Value = 114540

# Row 167
$ws.Range("G167").Value = 0.99036
$ws.Range("H167").Value = 0.3283
$ws.Range("I167").Value = 348
$ws.Range("J167").Value = 110600
$ws.Range("K167").Value = 1077
$ws.Range("L167").Value = 712
$ws.Range("M167").Value = 1971

# Row 168
$ws.Range("G168").Value = 0.27578
$ws.Range("H168").Value = 0.62245
$ws.Range("I168").Value = 122
$ws.Range("J168").Value = 1230
$ws.Range("K168").Value = 3230
$ws.Range("L168").Value = 74
$ws.Range("M168").Value = 110050

# Row 169
$ws.Range("G169").Value = 0.43457
$ws.Range("H169").Value = 0.51316
$ws.Range("I169").Value = 39
$ws.Range("J169").Value = 870
$ws.Range("K169").Value = 1132
$ws.Range("L169").Value = 37
$ws.Range("M169").Value = 112630

# Row 170
$ws.Range("G170").Value = 0.39432
$ws.Range("H170").Value = 0.79468
$ws.Range("I170").Value = 836
$ws.Range("J170").Value = 42894
$ws.Range("K170").Value = 65886
$ws.Range("L170").Value = 216
$ws.Range("M170").Value = 4872

# Row 171
$ws.Range("G171").Value = 0.39054
$ws.Range("H171").Value = 0.71003
$ws.Range("I171").Value = 524
$ws.Range("J171").Value = 12304
$ws.Range("K171").Value = 19201
$ws.Range("L171").Value = 214
$ws.Range("M171").Value = 82461

# Row 172
$ws.Range("G172").Value = 0.91155
$ws.Range("H172").Value = 0.32202
$ws.Range("I172").Value = 351
$ws.Range("J172").Value = 103520
$ws.Range("K172").Value = 10044
$ws.Range("L172").Value = 739
$ws.Range("M172").Value = 52

# Row 173
$ws.Range("G173").Value = 0.96949
$ws.Range("H173").Value = 0.71429
$ws.Range("I173").Value = 100
$ws.Range("J173").Value = 4926
$ws.Range("K173").Value = 155
$ws.Range("L173").Value = 40
$ws.Range("M173").Value = 109480

# Row 174
$ws.Range("G174").Value = 0.32022
$ws.Range("H174").Value = 0.81985
$ws.Range("I174").Value = 892
$ws.Range("J174").Value = 36322
$ws.Range("K174").Value = 77107
$ws.Range("L174").Value = 196
$ws.Range("M174").Value = 187

# Row 175
$ws.Range("G175").Value = 0.92473
$ws.Range("H175").Value = 0.65161
$ws.Range("I175").Value = 707
$ws.Range("J175").Value = 104650
$ws.Range("K175").Value = 8519
$ws.Range("L175").Value = 378
$ws.Range("M175").Value = 447

# Row 176
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 1
$ws.Range("I176").Value = 1091
$ws.Range("J176").Value = 0
$ws.Range("K176").Value = 113090
$ws.Range("L176").Value = 0
$ws.Range("M176").Value = 526

# Row 177
$ws.Range("G177").Value = 0.99908
$ws.Range("H177").Value = 0.00091996
$ws.Range("I177").Value = 1
$ws.Range("J177").Value = 112470
$ws.Range("K177").Value = 104
$ws.Range("L177").Value = 1086
$ws.Range("M177").Value = 1040

# Row 178
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 1
$ws.Range("I178").Value = 1091
$ws.Range("J178").Value = 0
$ws.Range("K178").Value = 113340
$ws.Range("L178").Value = 0
$ws.Range("M178").Value = 275

# Row 179
$ws.Range("G179").Value = 0.99987
$ws.Range("H179").Value = 0
$ws.Range("I179").Value = 0
$ws.Range("J179").Value = 112700
$ws.Range("K179").Value = 15
$ws.Range("L179").Value = 1087
$ws.Range("M179").Value = 907

# Row 180
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = "NaN"
$ws.Range("I180").Value = 0
$ws.Range("J180").Value = 0
$ws.Range("K180").Value = 57
$ws.Range("L180").Value = 0
$ws.Range("M180").Value = 114650

# Row 181
$ws.Range("G181").Value = 0.97746
$ws.Range("H181").Value = 0.417
$ws.Range("I181").Value = 309
$ws.Range("J181").Value = 43807
$ws.Range("K181").Value = 1010
$ws.Range("L181").Value = 432
$ws.Range("M181").Value = 69146

# Row 182
$ws.Range("G182").Value = 0.88178
$ws.Range("H182").Value = 0.56792
$ws.Range("I182").Value = 602
$ws.Range("J182").Value = 99930
$ws.Range("K182").Value = 13398
$ws.Range("L182").Value = 458
$ws.Range("M182").Value = 316

# Row 183
$ws.Range("G183").Value = 0.93928
$ws.Range("H183").Value = 0.59961
$ws.Range("I183").Value = 614
$ws.Range("J183").Value = 105170
$ws.Range("K183").Value = 6799
$ws.Range("L183").Value = 410
$ws.Range("M183").Value = 1708

# Row 184
$ws.Range("G184").Value = 0.0003363
$ws.Range("H184").Value = 0.99908
$ws.Range("I184").Value = 1090
$ws.Range("J184").Value = 38
$ws.Range("K184").Value = 112960
$ws.Range("L184").Value = 1
$ws.Range("M184").Value = 618

# Row 185
$ws.Range("G185").Value = 0.14145
$ws.Range("H185").Value = 0.68816
$ws.Range("I185").Value = 715
$ws.Range("J185").Value = 15871
$ws.Range("K185").Value = 96330
$ws.Range("L185").Value = 324
$ws.Range("M185").Value = 1464

# Row 186
$ws.Range("G186").Value = 0.95605
$ws.Range("H186").Value = 0.63477
$ws.Range("I186").Value = 690
$ws.Range("J186").Value = 108380
$ws.Range("K186").Value = 4982
$ws.Range("L186").Value = 397
$ws.Range("M186").Value = 253

# Row 187
$ws.Range("G187").Value = 0.94371
$ws.Range("H187").Value = 0.66269
$ws.Range("I187").Value = 723
$ws.Range("J187").Value = 106990
$ws.Range("K187").Value = 6382
$ws.Range("L187").Value = 368
$ws.Range("M187").Value = 237


# ---------------------------------------------------------------------------
# View state: move the visible/selected window down to the newly filled
# rows on both sheets (mirrors the author scrolling to see the new data).
# ---------------------------------------------------------------------------

# "Single Classification" (active sheet): scroll + select the newly
# populated block, G157:M187.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 141
$win.ScrollColumn = 1
$ws.Range("G157:M187").Select()

# "Group Classification": its stored selection also tracks the same row
# range (G157:M187); re-activate "Single Classification" afterwards so it
# stays the visible/active tab.
$ws2 = $wb.Worksheets.Item("Group Classification")
$ws2.Range("G157:M187").Select()
$ws.Activate()
